$d = $word.ActiveDocument

# --- Change 1: insert a new "Opportunity number: XXXX-XXXX-XXXX" paragraph
#     right after the "OpDiv: ..." paragraph, and before "Opportunity name: ..."
$opDivPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "OpDiv:*") {
        $opDivPara = $p
        break
    }
}
$opDivPara.Range.InsertParagraphAfter()

$newParaStart = $opDivPara.Range.End
$r1 = $d.Range($newParaStart, $newParaStart)
$r1.InsertAfter("Opportunity number: ")

$pos2 = $newParaStart + ("Opportunity number: ").Length
$r2 = $d.Range($pos2, $pos2)
$r2.InsertAfter("XXXX-XXXX-XXXX")

# Nudge formatting on just the number so it is kept in its own run
# (rather than being silently re-coalesced into the label's run).
$r3 = $d.Range($pos2, $pos2 + ("XXXX-XXXX-XXXX").Length)
$r3.Font.Bold = 1
$r3.Font.Bold = 0

# --- Change 2: collapse the split "Opportunity name" runs (", t" + "esting
#     fixture file") back into a single run with the same visible text.
$fix = $d.Content
$fix.Find.ClearFormatting()
$found = $fix.Find.Execute(", testing fixture file", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $fix.Text = ""
    $fix.InsertAfter(", testing fixture file")
}
